$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Test case 4 renumbered to 1 ---
$ws.Range("A20").Value = 1

# --- Remove existing hyperlinks before shifting rows (their refs won't auto-shift) ---
$ws.Hyperlinks.Delete()

# --- Insert a new row at 21: pushes old rows 21-25 down to 22-26 ---
$ws.Rows.Item(21).Insert()
$ws.Range("C21").Value = "age=22"
$ws.Range("D21").Value = "name=prasuna"
$ws.Range("E21").Value = "name=prasuna"

# Row 22 (old row 21 content shifted down): update D/E to age=22, add F22 = PASS
$ws.Range("D22").Value = "age=22"
$ws.Range("E22").Value = "age=22"
$ws.Range("F22").Value = "PASS"

# Row 23 (old row 22 content shifted down): remove stray F23 value
$ws.Range("F23").ClearContents()

# Row 25 (old row 24 content shifted down): add C25 = age=33
$ws.Range("C25").Value = "age=33"

# --- Insert another new row at 27 (after old row 25, now row 26) ---
$ws.Rows.Item(27).Insert()
$ws.Range("D27:E27").Style = "Normal"
$ws.Range("D27").Value = "age=33"
$ws.Range("E27").Value = "age=33"

# --- New test case 2 "help command" at row 31 ---
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "help command"
$ws.Range("C31").Value = "argv[0] -h"
$ws.Range("D31").Value = "enter inputs"
$ws.Range("E31").Value = "enter inputs"
$ws.Range("F31").Value = "PASS"

# --- Re-add hyperlinks at their new (shifted) locations, in the same rId order as before ---
$ws.Hyperlinks.Add($ws.Range("C22"), "mailto:mail=@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C24"), "mailto:mail=@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E23"), "mailto:mail=@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E26"), "mailto:mail=@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D23"), "mailto:mail=@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D26"), "mailto:mail=@gmail.com")

# Restore the shared "Hyperlink" cell style (Add() above mints a duplicate style variant)
$ws.Range("C22").Style = "Hyperlink"
$ws.Range("C24").Style = "Hyperlink"
$ws.Range("E23").Style = "Hyperlink"
$ws.Range("E26").Style = "Hyperlink"
$ws.Range("D23").Style = "Hyperlink"
$ws.Range("D26").Style = "Hyperlink"

# --- View state: scrolled down, new selection ---
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("B34").Select()
